$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8380
$ws.Range("I64").Value = 3900
$ws.Range("J64").Value = 9500
$ws.Range("K64").Value = 3900
$ws.Range("L64").Value = 9500
$ws.Range("M64").Value = -3652
$ws.Range("N64").Value = -9996

$ws.Range("H67").Value = 8380
$ws.Range("I67").Value = 3900
$ws.Range("J67").Value = 9500
$ws.Range("K67").Value = 3900
$ws.Range("L67").Value = 9500
$ws.Range("M67").Value = -3042
$ws.Range("N67").Value = -11216

$ws.Range("H74").Value = 3242.8333
$ws.Range("I74").Value = 3242.8333
$ws.Range("K74").Value = 3242.8333
$ws.Range("M74").Value = -2306.8333

$ws.Range("H77").Value = 3242.8333
$ws.Range("I77").Value = 3242.8333
$ws.Range("K77").Value = 16214.1665
$ws.Range("M77").Value = -11534.1665

$ws.Range("H86").Value = 5498.5
$ws.Range("J86").Value = 5498.5
$ws.Range("L86").Value = 5498.5
$ws.Range("N86").Value = -7744.5

$ws.Range("H89").Value = 5498.5
$ws.Range("J89").Value = 5498.5
$ws.Range("L89").Value = 27492.5
$ws.Range("N89").Value = -38724.5

$ws.Range("H106").Value = 7799.1
$ws.Range("I106").Value = 7553.8887
$ws.Range("K106").Value = 7553.8887
$ws.Range("M106").Value = -6922.8887

$ws.Range("H107").Value = 804
$ws.Range("I107").Value = 936.5
$ws.Range("J107").Value = 539
$ws.Range("K107").Value = 936.5
$ws.Range("L107").Value = 539
$ws.Range("M107").Value = 983.5
$ws.Range("N107").Value = -4379

$ws.Range("H115").Value = 892.5
$ws.Range("I115").Value = 892.5
$ws.Range("K115").Value = 2677.5
$ws.Range("M115").Value = -1110.5

$ws.Range("H127").Value = 1005.8333
$ws.Range("J127").Value = 495
$ws.Range("L127").Value = 1485
$ws.Range("N127").Value = -11405

$ws.Range("H129").Value = 1395
$ws.Range("I129").Value = 1198.5
$ws.Range("J129").Value = 1591.5
$ws.Range("K129").Value = 3595.5
$ws.Range("L129").Value = 4774.5
$ws.Range("M129").Value = 1404.5
$ws.Range("N129").Value = -14774.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 90
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -314

$ws.Range("H45").Value = 2701.4167
$ws.Range("I45").Value = 1760.75
$ws.Range("K45").Value = 1760.75
$ws.Range("M45").Value = -1383.75

$ws.Range("H61").Value = 4263.4
$ws.Range("I61").Value = 4263.4
$ws.Range("K61").Value = 4263.4
$ws.Range("M61").Value = -4051.4

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").ClearContents()

$ws.Range("H102").Value = 3897.611
$ws.Range("I102").Value = 1762.5834
$ws.Range("K102").Value = 1762.5834
$ws.Range("M102").Value = -140.5834

$ws.Range("H136").Value = 4263.4
$ws.Range("I136").Value = 4263.4
$ws.Range("K136").Value = 12790.2
$ws.Range("M136").Value = -10240.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 90
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -320

$ws.Range("H64").Value = 562.6
$ws.Range("I64").Value = 652
$ws.Range("K64").Value = 652
$ws.Range("M64").Value = -427

$ws.Range("H67").Value = 562.6
$ws.Range("I67").Value = 652
$ws.Range("K67").Value = 652
$ws.Range("M67").Value = 128

$ws.Range("H82").Value = 33354
$ws.Range("I82").Value = 14918.2
$ws.Range("K82").Value = 14918.2
$ws.Range("M82").Value = -14535.2

$ws.Range("H85").Value = 33354
$ws.Range("I85").Value = 14918.2
$ws.Range("K85").Value = 14918.2
$ws.Range("M85").Value = -13592.2

$ws.Range("H86").Value = 3781.96
$ws.Range("I86").Value = 2397.3572
$ws.Range("K86").Value = 2397.3572
$ws.Range("M86").Value = -1274.3572

$ws.Range("H89").Value = 3781.96
$ws.Range("I89").Value = 2397.3572
$ws.Range("K89").Value = 11986.786
$ws.Range("M89").Value = -6370.786

$ws.Range("H99").Value = 2622
$ws.Range("I99").Value = 1815.8
$ws.Range("K99").Value = 1815.8
$ws.Range("M99").Value = -317.8

$ws.Range("H107").Value = 4253.08
$ws.Range("I107").Value = 3197.4707
$ws.Range("K107").Value = 3197.4707
$ws.Range("M107").Value = -1277.4707

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 530
$ws.Range("I38").Value = 530
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 530
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

$ws.Range("H46").Value = 530
$ws.Range("I46").Value = 530
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 530
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H58").Value = 3892.2
$ws.Range("I58").Value = 1705.4286
$ws.Range("J58").Value = 8994.666999999999
$ws.Range("K58").Value = 1705.4286
$ws.Range("L58").Value = 8994.666999999999
$ws.Range("M58").Value = -1502.4286
$ws.Range("N58").Value = -9400.666999999999

$ws.Range("H134").Value = 1471
$ws.Range("I134").Value = 1651.6666
$ws.Range("K134").Value = 4954.9998
$ws.Range("M134").Value = -2419.9998

$ws.Range("H136").Value = 3892.2
$ws.Range("I136").Value = 1705.4286
$ws.Range("J136").Value = 8994.666999999999
$ws.Range("K136").Value = 5116.2858
$ws.Range("L136").Value = 26984.001
$ws.Range("M136").Value = -2566.2858
$ws.Range("N136").Value = -32084.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11899.5
$ws.Range("I3").Value = 11899.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 35698.5
$ws.Range("L3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("M3").Value = -35586.5

$ws.Range("H12").Value = 158.82353
$ws.Range("J12").Value = 175.41667
$ws.Range("L12").Value = 526.25001
$ws.Range("N12").Value = -872.25001

$ws.Range("H52").Value = 1174.5
$ws.Range("J52").Value = 1174.5
$ws.Range("L52").Value = 3523.5
$ws.Range("N52").Value = -4055.5

$ws.Range("H54").Value = 2104.4
$ws.Range("I54").Value = 2227.111
$ws.Range("K54").Value = 6681.333
$ws.Range("M54").Value = -6122.333

$ws.Range("H112").Value = 3500
$ws.Range("J112").Value = 3500
$ws.Range("L112").Value = 10500
$ws.Range("N112").Value = -12716

$ws.Range("H114").Value = 745.8333
$ws.Range("I114").Value = 295
$ws.Range("K114").Value = 885
$ws.Range("M114").Value = 2369

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H80").Value = 3512.7
$ws.Range("I80").Value = 3387.7144
$ws.Range("K80").Value = 3387.7144
$ws.Range("M80").Value = -2389.7144

$ws.Range("H83").Value = 3512.7
$ws.Range("I83").Value = 3387.7144
$ws.Range("K83").Value = 16938.572
$ws.Range("M83").Value = -11946.572

$ws.Range("H97").Value = 450.47058
$ws.Range("J97").Value = 649.6
$ws.Range("L97").Value = 649.6
$ws.Range("N97").Value = -1641.6

$ws.Range("H132").Value = 1350
$ws.Range("I132").Value = 1350
$ws.Range("K132").Value = 4050
$ws.Range("M132").Value = -1520

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 4999.5
$ws.Range("I18").Value = 4999.5
$ws.Range("K18").Value = 4999.5
$ws.Range("M18").Value = -4827.5

$ws.Range("H61").Value = 7855.4
$ws.Range("J61").Value = 7855.4
$ws.Range("L61").Value = 7855.4
$ws.Range("N61").Value = -8259.4

$ws.Range("H93").Value = 1046.6428
$ws.Range("I93").Value = 1105.2
$ws.Range("J93").Value = 900.25
$ws.Range("K93").Value = 1105.2
$ws.Range("L93").Value = 900.25
$ws.Range("M93").Value = 142.8
$ws.Range("N93").Value = -3396.25

$ws.Range("H113").Value = 7855.4
$ws.Range("J113").Value = 7855.4
$ws.Range("L113").Value = 7855.4
$ws.Range("N113").Value = -12195.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1379.6
$ws.Range("I96").Value = 1349.5
$ws.Range("K96").Value = 1349.5
$ws.Range("M96").Value = 23.5

$ws.Range("H107").Value = 562.8182
$ws.Range("I107").Value = 589.1
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 1767.3
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 152.6999999999998
$ws.Range("N107").Value = -4740

$ws.Range("H113").Value = 734.9231
$ws.Range("I113").Value = 482
$ws.Range("K113").Value = 1446
$ws.Range("M113").Value = 724
